$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: PAN number -> numeric value, holding -> 1999
$ws.Range("A2").Value = 4651561515
$ws.Range("B2").Value = 1999

# Update row 3: PAN number -> numeric value, holding -> 2999
$ws.Range("A3").Value = 1234455668
$ws.Range("B3").Value = 2999

# Remove row 4 entirely (clear contents) - data now only spans A1:B3
$ws.Range("A4:B4").ClearContents()

# Update the selection to match the new active cell
$ws.Range("B3").Select()
